# "Generate Report for Handback" -- refresh the localization-status report
# after a handback run: the files are back in sync with en-US, a fresh
# handback timestamp/file is recorded per-locale, and the stale
# "handback file is not latest" warning is cleared now that it is current.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Status column widened to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-03 00:53:30"
$zhcn.Range("P2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(16).ColumnWidth = 12.8333

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-09-03 00:53:38"
$dede.Range("P2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(16).ColumnWidth = 12.8333
